$d = $word.ActiveDocument

# The last paragraph currently reads "Crear escena de Game Over" and carries
# the hidden _GoBack bookmark right after "Over". We need to:
#   1. Leave that paragraph as-is (just losing the bookmark).
#   2. Add a new empty paragraph after it (same run formatting).
#   3. Add another new paragraph after that with the text
#      "Terminar de Mapear los Controles".
#   4. Move the _GoBack bookmark to the very end of that new last paragraph.

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range.Duplicate
$r.Collapse(0)

# Insert the first (empty) new paragraph.
$r.InsertParagraphAfter()

# Move to the newly created empty paragraph, then append the second new
# paragraph after it.
$emptyPara = $d.Paragraphs.Last
$r2 = $emptyPara.Range.Duplicate
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# The new last paragraph is where the heading text goes.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range.Duplicate
$newRange.Collapse(0)
$newRange.InsertAfter("Terminar de Mapear los Controles")

# Re-anchor the _GoBack bookmark at the very end of the new last paragraph.
# A bookmark collapsed exactly at a paragraph's trailing edge gets anchored
# to the paragraph start instead of staying collapsed, so insert a throwaway
# marker character, bookmark just before it, then remove the marker.
$endRange = $d.Paragraphs.Last.Range.Duplicate
$endRange.Collapse(0)
$endRange.InsertAfter("~")

$markerRange = $d.Content
$markerRange.Find.Execute("~") | Out-Null
$bmRange = $d.Range($markerRange.Start, $markerRange.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange2 = $d.Content
$markerRange2.Find.Execute("~") | Out-Null
$markerRange2.Text = ""
